$d = $word.ActiveDocument

# Helper: locate the start offset of the first occurrence of $text
function Get-StartPos([string]$text) {
    $r = $d.Content
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: text not found (start):" $text
    }
    return $r.Start
}

# Helper: locate the end offset of the first occurrence of $text
function Get-EndPos([string]$text) {
    $r = $d.Content
    $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "WARNING: text not found (end):" $text
    }
    return $r.End
}

# Add a bookmark named $name spanning [$startPos, $endPos)
function Add-Bookmark([string]$name, [int]$startPos, [int]$endPos) {
    $r = $d.Range($startPos, $endPos)
    $d.Bookmarks.Add($name, $r) | Out-Null
}

$ediaeresis = [string][char]0x00EB
$idiaeresis = [string][char]0x00EF

# Bookmark 0: _Hlk153195725 - "Wat doet de ecb?" heading through its answer paragraph
$s = Get-StartPos "Wat doet de ecb?"
$e = Get-EndPos "rentestanden het belangrijkste is."
Add-Bookmark "_Hlk153195725" $s $e

# Bookmark 1: _Hlk153195729 - "Waarom zijn die rentestanden zo belangrijk?" heading through its answer
$s = Get-StartPos "Waarom zijn die rentestanden zo belangrijk?"
$e = Get-EndPos "af willen remmen."
Add-Bookmark "_Hlk153195729" $s $e

# Bookmark 2: _Hlk153195734 - "2 jaar achter de rug..." through end of that paragraph
$s = Get-StartPos "2 jaar achter de rug van buitensporig"
$e = Get-EndPos "Dus ze zijn zo sowieso al verantwoordelijk. "
Add-Bookmark "_Hlk153195734" $s $e

# Bookmark 3: _Hlk153195744 - "Ja ze hadden het deels kunnen voorkomen..." through "...het tekort aan chips,"
$s = Get-StartPos "Ja ze hadden het deels kunnen voorkomen"
$e = Get-EndPos " helft van 2021, het tekort aan chips,"
Add-Bookmark "_Hlk153195744" $s $e

# Bookmark 4: _Hlk153195749 - "niet veroorzaakt maar wel gefaciliteerd..." through the space after "...prijsstijgingen de hebben."
$s = Get-StartPos "niet veroorzaakt maar wel gefaciliteerd"
$e = Get-EndPos "prijsstijgingen de hebben. "
Add-Bookmark "_Hlk153195749" $s $e

# Bookmark 5: _Hlk153195755 - "ECB heeft het zelf gefaciliteerd..." through "...die in omloop kwam"
$s = Get-StartPos "ECB heeft het zelf gefaciliteerd"
$e = Get-EndPos " die in omloop kwam"
Add-Bookmark "_Hlk153195755" $s $e

# Bookmark 6: _Hlk153195762 - "Maar daarnaast heb je ook nog ns een keertje..." run
$s = Get-StartPos "Maar daarnaast heb je ook nog ns een keertje"
$e = Get-EndPos "situatie gebruik maken door de prijzen nog wat extra te verhogen"
Add-Bookmark "_Hlk153195762" $s $e

# Bookmark 7: _Hlk153195769 - "Waardoor ... hele hoge economische groei gehad..." through "... om de prijzen extra te verhogen" (excludes trailing ".")
$s = Get-EndPos "huishouden gaan besteden. "
$e = Get-EndPos " om de prijzen extra te verhogen"
Add-Bookmark "_Hlk153195769" $s $e

# Bookmark 8: _Hlk153195691 - "In 2025 wil ecb..." heading through its answer paragraph
$s = Get-StartPos "In 2025 wil "
$e = Get-EndPos "dat dat wel gehaald gaat worden."
Add-Bookmark "_Hlk153195691" $s $e

# Bookmark 9: _Hlk153195700 - "Dus we krijgen volgend jaar eigenlijk pas..." through "... gedaan heeft" (excludes trailing ".")
$s = Get-StartPos "Dus we krijgen volgend jaar eigenlijk pas"
$e = Get-EndPos " gedaan heeft"
Add-Bookmark "_Hlk153195700" $s $e

# Bookmark 10: _Hlk153195705 - "Zijn stelling is dat de ECB reputatie verlies..." through "...geloofwaardigheid was flink ingezakt" (excludes trailing ".")
$s = Get-StartPos "Zijn stelling is dat de ECB reputatie verlies"
$e = Get-EndPos "geloofwaardigheid was flink ingezakt"
Add-Bookmark "_Hlk153195705" $s $e

# Bookmark 11: _Hlk153195710 - "Het inzakken van de economische groei..." through the space after "...we zitten in een crisis."
$s = Get-StartPos "Het inzakken van de economische groei"
$e = Get-EndPos "we zitten in een crisis. "
Add-Bookmark "_Hlk153195710" $s $e

# Bookmark 12: _Hlk153195715 - "Moeten we gaan leven met deze nieuwe prijs voor producten?" through end of its answer paragraph
$s = Get-StartPos "Moeten we gaan leven met deze nieuwe prijs voor producten?"
$e = Get-EndPos ("de inval van Oekra" + $idiaeresis + "ne.")
Add-Bookmark "_Hlk153195715" $s $e

Write-Host "Done adding bookmarks."
